# Applies the edits described by the commit:
#  1. Retitle the first proposal's "Project/Proposal Title" text.
#  2. Change the first proposal's summer person-months from "2" to "1.5".
#  3. Prefix the second proposal's title with "Collaborative Research: "
#     and move the "_GoBack" bookmark there (it used to sit in the
#     trailing empty paragraph at the very end of the document).

$d = $word.ActiveDocument

# --- 1. Update the first proposal's title text -----------------------
$rng = $d.Content
$null = $rng.Find.Execute(
    "Identifying the Influence of Anthropogenic Forcing on Extreme Weather Events",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "A hierarchical approach to improving the representation of convection in next-generation global models",
    2)

# --- 2. Update the first proposal's "Sumr:" person-months cell -------
# There are two identical "Sumr:  2" cells (one per proposal); target the
# one that belongs to the first proposal (table row 18, column 9).
$t = $d.Tables(1)
$sumrCell = $t.Cell(18, 9)
$sumrRange = $sumrCell.Range
$null = $sumrRange.Find.Execute(
    "2", $false, $false, $false, $false, $false, $true, 1, $false,
    "1.5", 2)

# --- 3. Prefix the second proposal's title and relocate the bookmark -
# Remove the bookmark from its old spot (trailing empty paragraph at the
# end of the document) before re-adding it elsewhere.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$titleRng = $d.Content
$null = $titleRng.Find.Execute(
    "Projecting future changes to lifecycle characteristics of North Atlantic cyclones")
$titleRng.Collapse(1)
$titleRng.InsertBefore("Collaborative Research: ")

# Re-find the (now shifted) title text so the bookmark lands immediately
# before it, matching the edited document.
$bmRng = $d.Content
$null = $bmRng.Find.Execute(
    "Projecting future changes to lifecycle characteristics of North Atlantic cyclones")
$bmRng.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRng)
